# fix mac_ids in apis
#
# The bulk-upload template used to key rows off the BLE MAC address; the
# API now wants the Keepr device id instead, and `product_uuid` becomes
# "<device_id> <uuid>". Update the sample/header data accordingly, swap
# the column order/widths (device_id first, narrower; product_uuid
# second, wider) and keep the rest of the blank templated rows below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Column widths: A = device_id (narrow), B = product_uuid (wide)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(2).ColumnWidth = 31.75

# ---------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "device_id"
$ws.Cells.Item(1,2).Value = "product_uuid"
$ws.Cells.Item(1,3).Value = "major"
$ws.Cells.Item(1,4).Value = "minor"
$ws.Cells.Item(1,5).Value = "color"

# ---------------------------------------------------------------------
# Sample data rows
# ---------------------------------------------------------------------
$ws.Cells.Item(2,1).Value = "KEEPR90628497"
$ws.Cells.Item(2,2).Value = "KEEPR90628497 EEF8EF65-AAAA-4410-B201-B6E1C4B9A486"
$ws.Cells.Item(2,3).Value = 9062
$ws.Cells.Item(2,4).Value = 8497
$ws.Cells.Item(2,5).Value = "White"

$ws.Cells.Item(3,1).Value = "KEEPR90638498"
$ws.Cells.Item(3,2).Value = "KEEPR90638498 EEF8EF65-AAAA-4410-B201-B6E1C4B9A486"
$ws.Cells.Item(3,3).Value = 9063
$ws.Cells.Item(3,4).Value = 8498
$ws.Cells.Item(3,5).Value = "Black"

# ---------------------------------------------------------------------
# Styling
# ---------------------------------------------------------------------
# Column A (device_id) gets the smaller Arial font used elsewhere in the
# template, left aligned. Build it once on a scratch cell and paste the
# whole format across so we don't bake extra throw-away styles in.
$scratch = $ws.Cells.Item(100, 26)
$scratch.Font.Size = 10
$scratch.Font.Name = "Arial"
$scratch.HorizontalAlignment = -4131
$scratch.Copy()
$ws.Range("A1:A3").PasteSpecial(-4122)
$scratch.Clear()
$excel.CutCopyMode = $false

# Header row + product_uuid column left aligned
$ws.Range("B1:E1").HorizontalAlignment = -4131
$ws.Range("B2:B3").HorizontalAlignment = -4131

# Numeric / value columns on the data rows right aligned
$ws.Range("C2:E3").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# Keep the remaining pre-formatted, blank template rows below the data
# ---------------------------------------------------------------------
for ($r = 4; $r -le 50; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

# Stray formatted (but empty) cell left over at D13
$ws.Cells.Item(13,4).WrapText = $false

# ---------------------------------------------------------------------
# Selection / view bookkeeping
# ---------------------------------------------------------------------
$ws.Range("B11").Select()
